$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "System Sequence Diagram" heading -> "Sequence Diagram"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "System Sequence Diagram", $true, $false, $false, $false, $false,
    $true, 1, $false, "Sequence Diagram", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Sequence-diagram hyperlink: new commit hash + new file name
# ---------------------------------------------------------------------------
$seqUrl = "https://github.com/aryansingh-ccm/BinaryBandits/blob/aca0310b3b280ef0c8853bde30e62040640ac31f/SystemDesign/LyricalLoomSequenceDiagram.vsdx"
$hSeq = $d.Hyperlinks(1)
$hSeq.TextToDisplay = $seqUrl
$hSeq.Address = $seqUrl

# ---------------------------------------------------------------------------
# 3) "Statechart Diagrams" heading -> "Statechart Diagram"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Statechart Diagrams", $true, $false, $false, $false, $false,
    $true, 1, $false, "Statechart Diagram", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Statechart-diagram hyperlink paragraph: new commit hash (file name keeps
#    the same), and the paragraph/run's redundant explicit 12pt size override
#    is dropped (it matched the document default anyway).  Rebuild the whole
#    paragraph via InsertXML so the stray w:sz/w:szCs are gone and the
#    hyperlink relationship is rewritten in one shot.
# ---------------------------------------------------------------------------
$stateUrl = "https://github.com/aryansingh-ccm/BinaryBandits/blob/aca0310b3b280ef0c8853bde30e62040640ac31f/SystemDesign/LyricalLoomStatechartDiagram.vsdx"

$statechartParaIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("https://github.com/aryansingh-ccm/BinaryBandits/blob/a80e4bb735e6daca570b3569661d2fe87b26f207/SystemDesign/LyricalLoomStatechartDiagram.vsdx")) {
        $statechartParaIndex = $i
    }
}
$pStatechart = $d.Paragraphs($statechartParaIndex)

$statechartXml = @"
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">
    <pkg:xmlData>
      <Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships>
    </pkg:xmlData>
  </pkg:part>
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
        <w:body>
          <w:p>
            <w:pPr>
              <w:rPr/>
            </w:pPr>
            <w:hyperlink r:id="rId2">
              <w:r>
                <w:rPr>
                  <w:color w:val="1155cc"/>
                  <w:u w:val="single"/>
                  <w:rtl w:val="0"/>
                </w:rPr>
                <w:t xml:space="preserve">$stateUrl</w:t>
              </w:r>
            </w:hyperlink>
            <w:r>
              <w:rPr>
                <w:rtl w:val="0"/>
              </w:rPr>
            </w:r>
          </w:p>
          <w:sectPr><w:pgSz w:w="12240" w:h="15840"/></w:sectPr>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
  <pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="256">
    <pkg:xmlData>
      <Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId2" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="$stateUrl" TargetMode="External"/></Relationships>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$pStatechart.Range.InsertXML($statechartXml) | Out-Null

# ---------------------------------------------------------------------------
# 5) Class-diagram hyperlink: new commit hash + renamed file
#    (LyricalLoomUpdatedClassDiagram.vsdx -> LyricalLoomClassDiagram.vsdx)
# ---------------------------------------------------------------------------
$classUrl = "https://github.com/aryansingh-ccm/BinaryBandits/blob/aca0310b3b280ef0c8853bde30e62040640ac31f/SystemDesign/LyricalLoomClassDiagram.vsdx"
$hClass = $d.Hyperlinks(3)
$hClass.TextToDisplay = $classUrl
$hClass.Address = $classUrl

# ---------------------------------------------------------------------------
# 6) Pseudocode hyperlink: new commit hash (file name unchanged)
# ---------------------------------------------------------------------------
$pseudoUrl = "https://github.com/aryansingh-ccm/BinaryBandits/blob/aca0310b3b280ef0c8853bde30e62040640ac31f/SystemDesign/LyricalLoomClassDiagramPseudocode.docx"
$hPseudo = $d.Hyperlinks(4)
$hPseudo.TextToDisplay = $pseudoUrl
$hPseudo.Address = $pseudoUrl

Write-Output "edit complete"
